$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "289÷6="
$t.Cell(1,2).Range.Text = "934÷2="
$t.Cell(1,3).Range.Text = "173÷8="
$t.Cell(1,4).Range.Text = "372÷9="
$t.Cell(1,5).Range.Text = "415÷4="
$t.Cell(5,1).Range.Text = "498÷5="
$t.Cell(5,2).Range.Text = "986÷3="
$t.Cell(5,3).Range.Text = "339÷6="
$t.Cell(5,4).Range.Text = "382÷8="
$t.Cell(5,5).Range.Text = "973÷8="
$t.Cell(9,1).Range.Text = "847÷3="
$t.Cell(9,2).Range.Text = "266÷9="
$t.Cell(9,3).Range.Text = "603÷6="
$t.Cell(9,4).Range.Text = "203÷4="
$t.Cell(9,5).Range.Text = "391÷7="
$t.Cell(13,1).Range.Text = "286÷4="
$t.Cell(13,2).Range.Text = "456÷6="
$t.Cell(13,3).Range.Text = "836÷5="
$t.Cell(13,4).Range.Text = "942÷5="
$t.Cell(13,5).Range.Text = "312÷8="
$t.Cell(17,1).Range.Text = "201÷4="
$t.Cell(17,2).Range.Text = "564÷4="
$t.Cell(17,3).Range.Text = "509÷2="
$t.Cell(17,4).Range.Text = "972÷8="
$t.Cell(17,5).Range.Text = "800÷4="
